# "selesai vr1 eps no" - add a "waktu sampai 500 games" timing column (H) to the
# analysis table, fill in the first run's High/Game-high scores and its elapsed
# time, bump the visual-range-1 row's D value, and bump the epsilon-reduction
# percentages from 1% to 2%.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the table with a new column H, inheriting the same borders/alignment
# already used by column G (header style + plain bordered body cells).
$ws.Range("G3:G17").Copy() | Out-Null
$ws.Range("H3:H17").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Fill in the new column's values (H4 first so the shared-string table mirrors
# the authoring order: "32.0m 6.95s" before "waktu sampai 500 games").
$ws.Range("H4").Value = "32.0m 6.95s"
$ws.Range("H3").Value = "waktu sampai 500 games"

# First run's High Score / Game high score results.
$ws.Range("F4").Value = 42
$ws.Range("G4").Value = 359

# Visual Range value for the second block (was 5, now 3).
$ws.Range("D6").Value = 3

# Epsilon reduction percentages bumped from 1% to 2% for every run.
$ws.Range("E5").Value = 0.02
$ws.Range("E7").Value = 0.02
$ws.Range("E9").Value = 0.02
$ws.Range("E11").Value = 0.02

# Widen column H and grow the (now taller, wrapped) header row to fit the new label.
$ws.Columns.Item(8).ColumnWidth = 16.833333333333336
$ws.Rows.Item(3).RowHeight = 45

$ws.Range("J8").Select() | Out-Null

Write-Host "edit applied"
